$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Refresh the "Elapsed Duration(Hrs)" column (G) on several sheets - the
#    report was regenerated later (~11h50m after the previous snapshot) so
#    every open-outage's elapsed duration grew by the same wall-clock gap.
# ---------------------------------------------------------------------------

$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3943:05:39"
$ws1.Range("G3").Value = "82:38:17"
$ws1.Range("G4").Value = "105:38:17"

$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12124:28:53"
$ws2.Range("G3").Value = "3254:12:22"
$ws2.Range("G4").Value = "492:23:56"

$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2970:18:42"
$ws4.Range("G3").Value = "197:30:57"
$ws4.Range("G4").Value = "85:43:22"
$ws4.Range("G5").Value = "83:20:55"

$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "444:17:41"

$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "84:49:59"

# ---------------------------------------------------------------------------
# 2. New outage row appended to sheet "R1" (row 6) - a fresh "NO PCM" style
#    entry for site JED0123, region R4, power source SCECO, status
#    "In progress", owner Latis (mirrors the blank-field pattern already
#    used by similar rows on this sheet).
# ---------------------------------------------------------------------------

$ws1.Range("B6").Value = "R4"
$ws1.Range("D6").Value = "JED0123"
$ws1.Range("I6").Value = "SCECO"
$ws1.Range("J6").Value = "In progress"
$ws1.Range("L6").Value = "Latis"
